# Auto-generated edit script applying the Golem_Profits.xlsx diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across
# the ALC, ARM, CRP, CUL, GSM, LTW, WVR sheets per the scheduled price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 599.8
$ws.Range("I2").Value = 433
$ws.Range("J2").Value = 850
$ws.Range("K2").Value = 433
$ws.Range("L2").Value = 850
$ws.Range("M2").Value = -320
$ws.Range("N2").Value = -1076
# Row 4
$ws.Range("H4").Value = 441.8
$ws.Range("I4").Value = 302.25
$ws.Range("K4").Value = 302.25
$ws.Range("M4").Value = -188.25
# Row 55
$ws.Range("H55").Value = 4564.5
$ws.Range("I55").Value = 4977.4
$ws.Range("K55").Value = 4977.4
$ws.Range("M55").Value = -4763.4
# Row 70
$ws.Range("H70").Value = 3111.111
$ws.Range("I70").Value = 3142.8572
$ws.Range("J70").Value = 3000
$ws.Range("K70").Value = 9428.571599999999
$ws.Range("L70").Value = 9000
$ws.Range("M70").Value = -9158.571599999999
$ws.Range("N70").Value = -9540
# Row 73
$ws.Range("H73").Value = 3111.111
$ws.Range("I73").Value = 3142.8572
$ws.Range("J73").Value = 3000
$ws.Range("K73").Value = 9428.571599999999
$ws.Range("L73").Value = 9000
$ws.Range("M73").Value = -8492.571599999999
$ws.Range("N73").Value = -10872
# Row 132
$ws.Range("H132").Value = 56588.168
$ws.Range("I132").Value = 56588.168
$ws.Range("K132").Value = 169764.504
$ws.Range("M132").Value = -167234.504

$ws = $wb.Worksheets.Item("ARM")
# Row 27
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
# Row 32
$ws.Range("H32").Value = 978
$ws.Range("I32").Value = 978
$ws.Range("K32").Value = 978
$ws.Range("M32").Value = -691
# Row 48
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
# Row 122
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
# Row 132
$ws.Range("H132").Value = 2008
$ws.Range("I132").Value = 1809.6
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 5428.799999999999
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -2898.799999999999
$ws.Range("N132").Value = -14060

$ws = $wb.Worksheets.Item("CRP")
# Row 17
$ws.Range("H17").Value = 250
$ws.Range("J17").Value = 250
$ws.Range("L17").Value = 250
$ws.Range("N17").Value = -598
# Row 43
$ws.Range("H43").Value = 11246
$ws.Range("J43").Value = 11246
$ws.Range("L43").Value = 11246
$ws.Range("N43").Value = -11614
# Row 101
$ws.Range("H101").Value = 11246
$ws.Range("J101").Value = 11246
$ws.Range("L101").Value = 11246
$ws.Range("N101").Value = -17736
# Row 134
$ws.Range("H134").Value = 1089
$ws.Range("I134").Value = 987.7778
$ws.Range("K134").Value = 2963.3334
$ws.Range("M134").Value = -428.3334
# Row 141
$ws.Range("H141").Value = 92141.57000000001
$ws.Range("J141").Value = 112998.4
$ws.Range("L141").Value = 112998.4
$ws.Range("N141").Value = -123358.4

$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
# Row 42
$ws.Range("H42").Value = 625
$ws.Range("J42").Value = 500
$ws.Range("L42").Value = 1500
$ws.Range("N42").Value = -2568
# Row 131
$ws.Range("H131").Value = 4197
$ws.Range("I131").Value = 1000
$ws.Range("J131").Value = 4996.25
$ws.Range("K131").Value = 3000
$ws.Range("L131").Value = 14988.75
$ws.Range("M131").Value = 2040
$ws.Range("N131").Value = -25068.75
# Row 132
$ws.Range("H132").Value = 690.8
$ws.Range("I132").Value = 574.75
$ws.Range("K132").Value = 5172.75
$ws.Range("M132").Value = -2642.75

$ws = $wb.Worksheets.Item("GSM")
# Row 112
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
# Row 122
$ws.Range("H122").Value = 7080.8335
$ws.Range("I122").Value = 5999.4
$ws.Range("J122").Value = 7853.2856
$ws.Range("K122").Value = 17998.2
$ws.Range("L122").Value = 23559.8568
$ws.Range("M122").Value = -15548.2
$ws.Range("N122").Value = -28459.8568

$ws = $wb.Worksheets.Item("LTW")
# Row 3
$ws.Range("H3").Value = 12000
$ws.Range("J3").Value = 12000
$ws.Range("L3").Value = 12000
$ws.Range("N3").Value = -12224
# Row 15
$ws.Range("H15").Value = 12000
$ws.Range("J15").Value = 12000
$ws.Range("L15").Value = 12000
$ws.Range("N15").Value = -12340
# Row 93
$ws.Range("H93").Value = 47620460
$ws.Range("J93").Value = 990.6667
$ws.Range("L93").Value = 990.6667
$ws.Range("N93").Value = -3486.6667
# Row 100
$ws.Range("H100").Value = 3346.6
$ws.Range("I100").Value = 3661.3333
$ws.Range("J100").Value = 2874.5
$ws.Range("K100").Value = 3661.3333
$ws.Range("L100").Value = 2874.5
$ws.Range("M100").Value = -3120.3333
$ws.Range("N100").Value = -3956.5

$ws = $wb.Worksheets.Item("WVR")
# Row 18
$ws.Range("H18").Value = 1000
$ws.Range("I18").Value = 1000
$ws.Range("K18").Value = 1000
$ws.Range("M18").Value = -827
# Row 26
$ws.Range("H26").Value = 948.5
$ws.Range("I26").Value = 833
$ws.Range("J26").Value = 1295
$ws.Range("K26").Value = 833
$ws.Range("L26").Value = 1295
$ws.Range("M26").Value = -540
$ws.Range("N26").Value = -1881
# Row 68
$ws.Range("H68").Value = 49332.668
$ws.Range("J68").Value = 49332.668
$ws.Range("L68").Value = 49332.668
$ws.Range("N68").Value = -50954.668
# Row 69
$ws.Range("H69").Value = 14249.556
$ws.Range("J69").Value = 14249.556
$ws.Range("L69").Value = 14249.556
$ws.Range("N69").Value = -15747.556
# Row 71
$ws.Range("H71").Value = 49332.668
$ws.Range("J71").Value = 49332.668
$ws.Range("L71").Value = 147998.004
$ws.Range("N71").Value = -156110.004
# Row 72
$ws.Range("H72").Value = 14249.556
$ws.Range("J72").Value = 14249.556
$ws.Range("L72").Value = 42748.66800000001
$ws.Range("N72").Value = -50236.66800000001
# Row 101
$ws.Range("H101").Value = 24801
$ws.Range("J101").Value = 24801
$ws.Range("L101").Value = 24801
$ws.Range("N101").Value = -31291
# Row 122
$ws.Range("H122").Value = 1979.8
$ws.Range("I122").Value = 1724.75
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 5174.25
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -2724.25
$ws.Range("N122").Value = -13900
